$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new rows for the Ueno shopping stops.
$ws.Range("A23").Value = "上野_多慶屋"
$ws.Range("A24").Value = "上野_yamashiroya玩具"
$ws.Range("D23").Value = "營業時間：11：00~23：00"
$ws.Range("D24").Value = "營業時間:10：00～21：30"

# Append a note about the 6th-floor observation deck to the existing
# "Tokyo station" direction/hours cell, and turn on wrap text so the
# added line displays on its own row.
$ws.Range("D16").Value = $ws.Range("D16").Value2 + " `n六樓有觀景台"
$ws.Range("D16").WrapText = $true
$ws.Rows(16).RowHeight = 33

# Leave the selection where the new content was added, matching the
# cursor position recorded by Excel on save.
$ws.Range("D21:D22").Select()
